$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The target run's original text is:
#   "144 Prot. n. 521963 del 19 dicembre 2025, in vigore dal 1° "
# It must become three runs (same rPr: rFonts eastAsia="Calibri"
# cstheme="minorHAnsi", szCs=20) with the protocol number corrected from
# 521963 to 531963:
#   "144 Prot. n. 5"  +  "3"  +  "1963 del 19 dicembre 2025, in vigore dal 1° "
# ---------------------------------------------------------------------------

$before  = "144 Prot. n. 5"
$digit   = "2"
$newDigit = "3"
$after   = "1963 del 19 dicembre 2025, in vigore dal 1"

# Locate the whole run by its (ASCII-only, degree-sign-free) text so the
# Find engine matches reliably.
$needle = $before + $digit + $after
$rng = $d.Content
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "target text not found" }

$runStart = $rng.Start
# Full original run length (includes the trailing '\u00b0 ' after "...dal 1").
$fullRunLength = $before.Length + $digit.Length + $after.Length + 2
$runEnd = $runStart + $fullRunLength

# Split points inside/around the run we are rewriting.
$splitA = $runStart + $before.Length                 # between "...n. 5" and the digit
$splitB = $splitA + $digit.Length                     # between the digit and "1963..."

# Replace the single digit "2" -> "3" in place (mid-run text edit).
$digitRange = $d.Range($splitA, $splitB)
$digitRange.Text = $newDigit

# The text edit above causes the whole paragraph run-chain (runs that share
# identical formatting) to coalesce into a single run. Re-impose the
# original run boundaries - both the new 3-way split of the edited run and
# the pre-existing boundaries of the runs that follow it ("febbraio", " ",
# "2026", ";") - by toggling a direct character-formatting property on and
# back off over each sub-range. Applying this "seal" from the rightmost
# boundary back to the leftmost keeps earlier seals from being undone by
# later ones.

function Seal-Boundary($rStart, $rEnd) {
    $r = $d.Range($rStart, $rEnd)
    $r.Bold = 1
    $r.Bold = 0
}

# Boundaries of the runs following the edited one: "febbraio", " ", "2026", ";"
$febWord  = "febbraio"
$space    = " "
$year     = "2026"
$semi     = ";"

$febStart  = $runEnd
$febEnd    = $febStart + $febWord.Length
$spaceEnd  = $febEnd + $space.Length
$yearEnd   = $spaceEnd + $year.Length
$semiEnd   = $yearEnd + $semi.Length

# Seal from right to left.
Seal-Boundary $yearEnd  $semiEnd    # isolates ";"
Seal-Boundary $spaceEnd $yearEnd    # isolates "2026"
Seal-Boundary $febEnd   $spaceEnd   # isolates " "
Seal-Boundary $runEnd   $febEnd     # isolates "febbraio"
Seal-Boundary $splitB   $runEnd     # isolates "1963 del 19 dicembre 2025, in vigore dal 1\u00b0 "
Seal-Boundary $splitA   $splitB     # isolates "3"
Seal-Boundary $runStart $splitA     # isolates "144 Prot. n. 5"

Write-Output "Protocol number corrected and run split applied."
